# Applies the "combat-masters (Version 2)" edit:
#   1. Remove the standalone "Meta description: ..." paragraph that
#      followed the H1 title.
#   2. Before the trailing (italic) image-prompt paragraph, insert a new
#      bold paragraph repeating the page title
#      "Play Combat Masters Slot Free - Review & Gameplay Features".
#   3. Replace the trailing paragraph's text (the old AI image prompt)
#      with the meta-description copy (minus the "Meta description: "
#      label), keeping its italic formatting.

$d = $word.ActiveDocument

# --- Step 1: delete the "Meta description" paragraph (2nd paragraph) ---
$metaLabel = "Meta description"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.StartsWith($metaLabel)) {
        $para.Range.Delete() | Out-Null
        break
    }
}

# --- Step 2: insert a new bold paragraph before the final paragraph ---
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$lastPara.Range.InsertParagraphBefore() | Out-Null

$titleParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Combat Masters Slot Free - Review &amp; Gameplay Features</w:t></w:r></w:p>'
$d.Paragraphs($count).Range.InsertXML($titleParaXml) | Out-Null

# --- Step 3: replace the final (italic) paragraph's text ---
$finalPara = $d.Paragraphs($d.Paragraphs.Count)
$descriptionXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Combat Masters with pros &amp; cons, gameplay features, payouts &amp; RTP. Play Combat Masters slot for free.</w:t></w:r></w:p>'
$finalPara.Range.InsertXML($descriptionXml) | Out-Null
